$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 34
$ws.Range("H34").Value = 20434
$ws.Range("I34").Value = 15649.75
$ws.Range("K34").Value = 15649.75
$ws.Range("M34").Value = -15446.75

# Row 36
$ws.Range("H36").Value = 20434
$ws.Range("I36").Value = 15649.75
$ws.Range("K36").Value = 15649.75
$ws.Range("M36").Value = -14934.75

# Row 53
$ws.Range("H53").Value = 353.0909
$ws.Range("I53").Value = 394
$ws.Range("J53").Value = 304
$ws.Range("K53").Value = 394
$ws.Range("L53").Value = 304
$ws.Range("M53").Value = 243
$ws.Range("N53").Value = -1578

# Row 61
$ws.Range("H61").Value = 550.6667
$ws.Range("I61").Value = 550.6667
$ws.Range("K61").Value = 1652.0001
$ws.Range("M61").Value = -1480.0001


# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 36
$ws.Range("H36").Value = 3502.8572
$ws.Range("I36").Value = 3502.8572
$ws.Range("K36").Value = 3502.8572
$ws.Range("M36").Value = -3156.8572

# Row 88
$ws.Range("H88").Value = 2760.182
$ws.Range("I88").Value = 2545.75
$ws.Range("K88").Value = 2545.75
$ws.Range("M88").Value = -2139.75

# Row 91
$ws.Range("H91").Value = 2760.182
$ws.Range("I91").Value = 2545.75
$ws.Range("K91").Value = 2545.75
$ws.Range("M91").Value = -1141.75

# Row 102
$ws.Range("H102").Value = 1267.75
$ws.Range("I102").Value = 1234.5714
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1234.5714
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 387.4286
$ws.Range("N102").Value = -4744

# Row 122
$ws.Range("H122").Value = 1968.25
$ws.Range("I122").Value = 1807.1538
$ws.Range("J122").Value = 2666.3333
$ws.Range("K122").Value = 5421.4614
$ws.Range("L122").Value = 7998.999899999999
$ws.Range("M122").Value = -2971.4614
$ws.Range("N122").Value = -12898.9999

# Row 132
$ws.Range("H132").Value = 1580.8214
$ws.Range("J132").Value = 2265.5715
$ws.Range("L132").Value = 6796.7145
$ws.Range("N132").Value = -11856.7145


# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 55
$ws.Range("H55").Value = 15000
$ws.Range("I55").Value = 15000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 15000
$ws.Range("M55").Value = -14727
$ws.Range("N55").ClearContents()

# Row 134
$ws.Range("H134").Value = 2588.2415
$ws.Range("I134").Value = 2847.16
$ws.Range("K134").Value = 8541.48
$ws.Range("M134").Value = -6006.48


# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 33
$ws.Range("H33").Value = 8938.888999999999
$ws.Range("I33").Value = 6311.2
$ws.Range("J33").Value = 12223.5
$ws.Range("K33").Value = 6311.2
$ws.Range("L33").Value = 12223.5
$ws.Range("M33").Value = -5932.2
$ws.Range("N33").Value = -12981.5

# Row 58
$ws.Range("H58").Value = 1956.8462
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

# Row 62
$ws.Range("H62").Value = 34857.684
$ws.Range("I62").Value = 3499.8462
$ws.Range("J62").Value = 102799.664
$ws.Range("K62").Value = 3499.8462
$ws.Range("L62").Value = 102799.664
$ws.Range("M62").Value = -2875.8462
$ws.Range("N62").Value = -104047.664

# Row 65
$ws.Range("H65").Value = 34857.684
$ws.Range("I65").Value = 3499.8462
$ws.Range("J65").Value = 102799.664
$ws.Range("K65").Value = 17499.231
$ws.Range("L65").Value = 513998.32
$ws.Range("M65").Value = -14379.231
$ws.Range("N65").Value = -520238.32

# Row 105
$ws.Range("H105").Value = 1949
$ws.Range("I105").Value = 1748.375
$ws.Range("K105").Value = 1748.375
$ws.Range("M105").Value = -1.375

# Row 134
$ws.Range("H134").Value = 2683.6428
$ws.Range("I134").Value = 2523.8147
$ws.Range("J134").Value = 6999
$ws.Range("K134").Value = 7571.4441
$ws.Range("L134").Value = 20997
$ws.Range("M134").Value = -5036.4441
$ws.Range("N134").Value = -26067

# Row 136
$ws.Range("H136").Value = 1956.8462
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()


# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 2508.6667
$ws.Range("J122").Value = 3851.6667
$ws.Range("L122").Value = 34665.0003
$ws.Range("N122").Value = -39565.0003

# Row 131
$ws.Range("H131").Value = 1407.4426
$ws.Range("I131").Value = 755.5833
$ws.Range("J131").Value = 1567.0817
$ws.Range("K131").Value = 2266.7499
$ws.Range("L131").Value = 4701.2451
$ws.Range("M131").Value = 2773.2501
$ws.Range("N131").Value = -14781.2451


# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 105
$ws.Range("H105").Value = 27098
$ws.Range("J105").Value = 27098
$ws.Range("L105").Value = 27098
$ws.Range("N105").Value = -34086

# Row 122
$ws.Range("H122").Value = 2653.0715
$ws.Range("I122").Value = 2724.4736
$ws.Range("J122").Value = 2502.3333
$ws.Range("K122").Value = 8173.4208
$ws.Range("L122").Value = 7506.999899999999
$ws.Range("M122").Value = -5723.4208
$ws.Range("N122").Value = -12406.9999

# Row 126
$ws.Range("H126").Value = 5006.0586
$ws.Range("I126").Value = 3086.1428
$ws.Range("J126").Value = 6350
$ws.Range("K126").Value = 9258.428400000001
$ws.Range("L126").Value = 19050
$ws.Range("M126").Value = -6788.428400000001
$ws.Range("N126").Value = -23990

# Row 132
$ws.Range("H132").Value = 1920.697
$ws.Range("I132").Value = 1775.4642
$ws.Range("K132").Value = 5326.392599999999
$ws.Range("M132").Value = -2796.392599999999


# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 149.42857
$ws.Range("I55").Value = 136
$ws.Range("K55").Value = 136
$ws.Range("M55").Value = 37

# Row 132
$ws.Range("H132").Value = 5711
$ws.Range("I132").Value = 4530.077
$ws.Range("J132").Value = 10828.333
$ws.Range("K132").Value = 13590.231
$ws.Range("L132").Value = 32484.999
$ws.Range("M132").Value = -11060.231
$ws.Range("N132").Value = -37544.999


# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 36327
$ws.Range("J75").Value = 36192.6
$ws.Range("L75").Value = 36192.6
$ws.Range("N75").Value = -38064.6

# Row 78
$ws.Range("H78").Value = 36327
$ws.Range("J78").Value = 36192.6
$ws.Range("L78").Value = 108577.8
$ws.Range("N78").Value = -117937.8

# Row 122
$ws.Range("H122").Value = 4130.4116
$ws.Range("I122").Value = 3998.4482
$ws.Range("J122").Value = 4895.8
$ws.Range("K122").Value = 11995.3446
$ws.Range("L122").Value = 14687.4
$ws.Range("M122").Value = -9545.3446
$ws.Range("N122").Value = -19587.4

# Row 132
$ws.Range("H132").Value = 3816.88
$ws.Range("I132").Value = 2444
$ws.Range("J132").Value = 4589.125
$ws.Range("K132").Value = 7332
$ws.Range("L132").Value = 13767.375
$ws.Range("M132").Value = -4802
$ws.Range("N132").Value = -18827.375

